$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for rows 194-212 (dates in column A, scheduled/tracked flights in B/C) ---
$newRows = @(
    @{ Row=194; Date="2020-10-16"; B=56; C=54 },
    @{ Row=195; Date="2020-10-17"; B=43; C=43 },
    @{ Row=196; Date="2020-10-18"; B=53; C=51 },
    @{ Row=197; Date="2020-10-19"; B=56; C=52 },
    @{ Row=198; Date="2020-10-20"; B=45; C=43 },
    @{ Row=199; Date="2020-10-21"; B=48; C=46 },
    @{ Row=200; Date="2020-10-22"; B=51; C=47 },
    @{ Row=201; Date="2020-10-23"; B=56; C=56 },
    @{ Row=202; Date="2020-10-24"; B=47; C=44 },
    @{ Row=203; Date="2020-10-25"; B=48; C=43 },
    @{ Row=204; Date="2020-10-26"; B=48; C=48 },
    @{ Row=205; Date="2020-10-27"; B=45; C=43 },
    @{ Row=206; Date="2020-10-28"; B=44; C=42 },
    @{ Row=207; Date="2020-10-29"; B=55; C=50 },
    @{ Row=208; Date="2020-10-30"; B=55; C=53 },
    @{ Row=209; Date="2020-10-31"; B=32; C=31 },
    @{ Row=210; Date="2020-11-01"; B=42; C=39 },
    @{ Row=211; Date="2020-11-02"; B=43; C=41 },
    @{ Row=212; Date="2020-11-03"; B=43; C=41 }
)

# Copy formatting from the last pre-existing data row (192) down across the new rows
# for columns A, B and C so the new cells reuse the existing styles instead of
# Excel inventing brand-new ones.
$ws.Range("A192:C192").Copy()
$ws.Range("A194:C212").PasteSpecial(-4122)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Date
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}

# --- Percentage formulas in column D ---
# Row 193 previously had no formula; it now continues the existing pattern,
# together with the freshly added rows 194-196.
$ws.Cells.Item(192, 4).Copy()
$ws.Range("D193:D196").PasteSpecial(-4122)
$ws.Range("D193:D196").Formula = "=C193/B193"

# Rows 197-211 form a new block of the same percentage formula.
$ws.Range("D197:D211").PasteSpecial(-4122)
$ws.Range("D197:D211").Formula = "=C197/B197"

# Row 212 has no percentage formula, matching the source data.

# --- Selection / view changes ---
# Scroll the sheet so row 180 is at the top of the visible window, and put
# the active cell/selection on H209 (matches the author's viewport at save time).
try {
    $excel.ActiveWindow.ScrollRow = 180
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("H209").Select()

# Move the workbook window itself (best-effort; matches the new xWindow/yWindow
# recorded in the workbook's bookViews).
try {
    $excel.ActiveWindow.Left = 35900
    $excel.ActiveWindow.Top = 1080
} catch {}
